$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-26 Thursday", 1) | Out-Null
$d.Content.Find.Execute("94×39=3666", $true, $false, $false, $false, $false, $true, 1, $false, "17×91=1547", 1) | Out-Null
$d.Content.Find.Execute("22×25=550", $true, $false, $false, $false, $false, $true, 1, $false, "40×47=1880", 1) | Out-Null
$d.Content.Find.Execute("30×16=480", $true, $false, $false, $false, $false, $true, 1, $false, "89×99=8811", 1) | Out-Null
$d.Content.Find.Execute("52×94=4888", $true, $false, $false, $false, $false, $true, 1, $false, "56×68=3808", 1) | Out-Null
$d.Content.Find.Execute("82×96=7872", $true, $false, $false, $false, $false, $true, 1, $false, "42×90=3780", 1) | Out-Null
$d.Content.Find.Execute("93×18=1674", $true, $false, $false, $false, $false, $true, 1, $false, "73×78=5694", 1) | Out-Null
$d.Content.Find.Execute("12×55=660", $true, $false, $false, $false, $false, $true, 1, $false, "32×96=3072", 1) | Out-Null
$d.Content.Find.Execute("63×66=4158", $true, $false, $false, $false, $false, $true, 1, $false, "14×85=1190", 1) | Out-Null
$d.Content.Find.Execute("95×86=8170", $true, $false, $false, $false, $false, $true, 1, $false, "54×35=1890", 1) | Out-Null
$d.Content.Find.Execute("98×16=1568", $true, $false, $false, $false, $false, $true, 1, $false, "93×15=1395", 1) | Out-Null
$d.Content.Find.Execute("36×37=1332", $true, $false, $false, $false, $false, $true, 1, $false, "95×76=7220", 1) | Out-Null
$d.Content.Find.Execute("84×93=7812", $true, $false, $false, $false, $false, $true, 1, $false, "47×65=3055", 1) | Out-Null
$d.Content.Find.Execute("31×95=2945", $true, $false, $false, $false, $false, $true, 1, $false, "32×77=2464", 1) | Out-Null
$d.Content.Find.Execute("70×73=5110", $true, $false, $false, $false, $false, $true, 1, $false, "60×21=1260", 1) | Out-Null
$d.Content.Find.Execute("15×71=1065", $true, $false, $false, $false, $false, $true, 1, $false, "71×30=2130", 1) | Out-Null
$d.Content.Find.Execute("48×26=1248", $true, $false, $false, $false, $false, $true, 1, $false, "12×55=660", 1) | Out-Null
$d.Content.Find.Execute("56×47=2632", $true, $false, $false, $false, $false, $true, 1, $false, "79×59=4661", 1) | Out-Null
$d.Content.Find.Execute("17×30=510", $true, $false, $false, $false, $false, $true, 1, $false, "68×11=748", 1) | Out-Null
$d.Content.Find.Execute("11×67=737", $true, $false, $false, $false, $false, $true, 1, $false, "80×36=2880", 1) | Out-Null
$d.Content.Find.Execute("54×30=1620", $true, $false, $false, $false, $false, $true, 1, $false, "51×89=4539", 1) | Out-Null
$d.Content.Find.Execute("71×20=1420", $true, $false, $false, $false, $false, $true, 1, $false, "38×82=3116", 1) | Out-Null
$d.Content.Find.Execute("81×21=1701", $true, $false, $false, $false, $false, $true, 1, $false, "56×38=2128", 1) | Out-Null
$d.Content.Find.Execute("58×54=3132", $true, $false, $false, $false, $false, $true, 1, $false, "61×76=4636", 1) | Out-Null
$d.Content.Find.Execute("50×57=2850", $true, $false, $false, $false, $false, $true, 1, $false, "76×52=3952", 1) | Out-Null
$d.Content.Find.Execute("83×65=5395", $true, $false, $false, $false, $false, $true, 1, $false, "59×17=1003", 1) | Out-Null
